# New weekly price record for Ajo (Chino / Primera) needs to be inserted
# ahead of the existing history, so every subsequent record shifts down
# by one row (old row 70 -> new row 71, ..., old row 131 -> new row 132).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(70).Insert()

$ws.Range("A70").Value = 7
$ws.Range("B70").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C70").Value = "Ñuble"
$ws.Range("D70").Value = 44484
$ws.Range("E70").Value = 16
$ws.Range("F70").Value = 100112003
$ws.Range("G70").Value = "Ajo"
$ws.Range("H70").Value = "Chino"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 60
$ws.Range("K70").Value = 16000
$ws.Range("L70").Value = 17000
$ws.Range("M70").Value = 16500
$ws.Range("N70").Value = "$/caja 10 kilos"
$ws.Range("O70").Value = "China"
$ws.Range("P70").Value = 1650
$ws.Range("Q70").Value = 10
$ws.Range("R70").Value = "Hortaliza"

$ws.Range("D70").NumberFormat = "YYYY-MM-DD HH:MM:SS"
